# Update parameter names (rename) in the "grwat" params workbook.
#
# Sheet "parameters" column C ("name") holds the short machine name of each
# parameter. Row 28 ("ModeMountain") still carried the old short name "mnt";
# rename it to the new name "mntmode".
#
# Sheet "values" row 1 is a header duplicating the short parameter names used
# as column headers; it still carried the OLD short names, so it needs to be
# refreshed to the current names (matching sheet "parameters" column C),
# including the same "mnt" -> "mntmode" rename.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("parameters")
$wsValues = $wb.Worksheets.Item("values")

# --- Rename parameter "mnt" -> "mntmode" on the parameters sheet ---
$wsParams.Range("C28").Value = "mntmode"

# --- Refresh the header row on the "values" sheet to the current parameter names ---
$wsValues.Range("B1").Value = "winmon"
$wsValues.Range("C1").Value = "grad1"
$wsValues.Range("D1").Value = "grad2"
$wsValues.Range("E1").Value = "gratio"
$wsValues.Range("F1").Value = "ftmon1"
$wsValues.Range("G1").Value = "ftmon2"
$wsValues.Range("H1").Value = "ftrisedays1"
$wsValues.Range("I1").Value = "ftrisedays2"
$wsValues.Range("J1").Value = "ftdays"
$wsValues.Range("K1").Value = "ftrise"
$wsValues.Range("L1").Value = "ftratio"
$wsValues.Range("M1").Value = "ftrecdays"
$wsValues.Range("N1").Value = "precdays"
$wsValues.Range("O1").Value = "frostdays"
$wsValues.Range("P1").Value = "windays"
$wsValues.Range("Q1").Value = "floodprec"
$wsValues.Range("R1").Value = "floodtemp"
$wsValues.Range("S1").Value = "frosttemp"
$wsValues.Range("T1").Value = "wintemp"
$wsValues.Range("U1").Value = "signratio1"
$wsValues.Range("V1").Value = "signratio2"
$wsValues.Range("W1").Value = "gapflag"
$wsValues.Range("X1").Value = "floodratio"
$wsValues.Range("Y1").Value = "gaplen"
$wsValues.Range("Z1").Value = "snowtemp"
$wsValues.Range("AB1").Value = "mntmode"
$wsValues.Range("AC1").Value = "mntgrad"
$wsValues.Range("AD1").Value = "mntavgdays"
$wsValues.Range("AE1").Value = "mntratiodays"
$wsValues.Range("AF1").Value = "mntratio"

# --- Update sheet view / selection state ---
# "parameters" is no longer the selected tab; its view is scrolled down a bit
# and the selection moved to C29.
$wsParams.Activate()
$wsParams.Application.ActiveWindow.ScrollRow = 2
$wsParams.Range("C29").Select()

# "values" becomes the active/selected tab, scrolled right, with AB1 selected.
$wsValues.Activate()
$wsValues.Application.ActiveWindow.ScrollColumn = 20
$wsValues.Range("AB1").Select()
